$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Test"

function Set-TextValue($row, $col, $val) {
  $c = $ws.Cells.Item($row, $col)
  $c.NumberFormat = "@"
  $c.Value = $val
  $c.ClearFormats()
}

# Row 1
$ws.Cells.Item(1,1).Value = "Building"
$ws.Cells.Item(1,2).Value = "Floor"
$ws.Cells.Item(1,3).Value = "Room"
$ws.Cells.Item(1,4).Value = "Storage"
$ws.Cells.Item(1,5).Value = "Barcode"
$ws.Cells.Item(1,6).Value = "Chemical Name"
$ws.Cells.Item(1,7).Value = "Current Quantity"
$ws.Cells.Item(1,8).Value = "Current Unit"

# Row 2
$ws.Cells.Item(2,1).Value = "Science"
$ws.Cells.Item(2,2).Value = "Basement"
$ws.Cells.Item(2,3).Value = "13A"
$ws.Cells.Item(2,4).Value = "13A"
$ws.Cells.Item(2,5).Value = "E1-3299"
$ws.Cells.Item(2,6).Value = "PHOSPHATE"
$ws.Cells.Item(2,7).Value = 22
$ws.Cells.Item(2,8).Value = "gram (g)"

# Row 3
$ws.Cells.Item(3,1).Value = "Science"
$ws.Cells.Item(3,2).Value = "Basement"
$ws.Cells.Item(3,3).Value = "13A"
$ws.Cells.Item(3,4).Value = "13A"
$ws.Cells.Item(3,5).Value = "N1-15"
$ws.Cells.Item(3,6).Value = "2-CHLOROBUTANE"
$ws.Cells.Item(3,7).Value = 22
$ws.Cells.Item(3,8).Value = "gram (g)"

# Row 4
$ws.Cells.Item(4,1).Value = "Science"
$ws.Cells.Item(4,2).Value = "Basement"
$ws.Cells.Item(4,3).Value = "13A"
$ws.Cells.Item(4,4).Value = "13A"
Set-TextValue 4 5 "17040000"
$ws.Cells.Item(4,6).Value = "2,6-DICHLOROINDOPHENOL SODIUM DERIVATIVE"
$ws.Cells.Item(4,7).Value = 33
$ws.Cells.Item(4,8).Value = "gram (g)"

# Row 5
$ws.Cells.Item(5,1).Value = "Science"
$ws.Cells.Item(5,2).Value = "Basement"
$ws.Cells.Item(5,3).Value = "13A"
$ws.Cells.Item(5,4).Value = "Cyanide Cabinet"
$ws.Cells.Item(5,5).Value = "E1-3399"
$ws.Cells.Item(5,6).Value = "2-(DIMETHYLAMINO)PYRIDINE"
$ws.Cells.Item(5,7).Value = 22
$ws.Cells.Item(5,8).Value = "gram (g)"

# Row 6
$ws.Cells.Item(6,1).Value = "Science"
$ws.Cells.Item(6,2).Value = "Basement"
$ws.Cells.Item(6,3).Value = "13C"
$ws.Cells.Item(6,4).Value = "13C"
$ws.Cells.Item(6,5).Value = "E1-3540"
$ws.Cells.Item(6,6).Value = "2-BROMOBUTANE"
$ws.Cells.Item(6,7).Value = 22
$ws.Cells.Item(6,8).Value = "gram (g)"

# Row 7
$ws.Cells.Item(7,1).Value = "Science"
$ws.Cells.Item(7,2).Value = "Basement"
$ws.Cells.Item(7,3).Value = "13C"
$ws.Cells.Item(7,4).Value = "13C"
Set-TextValue 7 5 "16020008"
$ws.Cells.Item(7,6).Value = "ACETIC ACID 99+%"
$ws.Cells.Item(7,7).Value = 33
$ws.Cells.Item(7,8).Value = "gram (g)"

# Row 8
$ws.Cells.Item(8,1).Value = "Science"
Set-TextValue 8 2 "2"
Set-TextValue 8 3 "203"
Set-TextValue 8 4 "203"
$ws.Cells.Item(8,5).Value = "D1-8364"
$ws.Cells.Item(8,6).Value = "2,6-DICHLOROINDOPHENOL SODIUM DERIVATIVE"
$ws.Cells.Item(8,7).Value = 22
$ws.Cells.Item(8,8).Value = "gram (g)"

Write-Host "done"
